# Update the "profit" value used across the workbook.
# CALCULADORA!B3 holds the numeric value (772.37 -> 771.17).
$wb = $excel.ActiveWorkbook

$wsCalc = $wb.Worksheets.Item("CALCULADORA")
$wsCalc.Range("B3").Value = 771.17

# "grilla de pruebas"!B3 stores the same number as text (shared string),
# referenced by several formulas on that sheet (e.g. -B3*$B$10/100).
# Update it so the sheet's dependent formulas recalculate consistently.
$wsGrilla = $wb.Worksheets.Item("grilla de pruebas")
$wsGrilla.Range("B3").Value = "771.17"

# Selection / active sheet bookkeeping: CALCULADORA becomes the active tab,
# "grilla de pruebas" keeps selection on B3 but is no longer the active tab.
$wsGrilla.Range("B3").Select()

$wsCalc.Activate()
$wsCalc.Range("B3").Select()
